# Update the "想去人数" (F column) values across the four worksheets to
# reflect the newly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1723
$ws1.Range("F4").Value  = 1252
$ws1.Range("F6").Value  = 154
$ws1.Range("F7").Value  = 1417
$ws1.Range("F9").Value  = 25
$ws1.Range("F10").Value = 667
$ws1.Range("F16").Value = 509
$ws1.Range("F21").Value = 737
$ws1.Range("F27").Value = 179
$ws1.Range("F33").Value = 110

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 635
$ws2.Range("F14").Value = 493

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value  = 2306
$ws3.Range("F9").Value  = 1155
$ws3.Range("F10").Value = 282

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 2306
$ws4.Range("F6").Value  = 1723
$ws4.Range("F9").Value  = 1155
$ws4.Range("F10").Value = 282
$ws4.Range("F12").Value = 1252
$ws4.Range("F14").Value = 154
$ws4.Range("F15").Value = 1417
$ws4.Range("F17").Value = 25
$ws4.Range("F18").Value = 667
$ws4.Range("F22").Value = 509
$ws4.Range("F26").Value = 737
$ws4.Range("F33").Value = 179
$ws4.Range("F38").Value = 493
$ws4.Range("F41").Value = 110
